$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark BOJ problems in column B as solved ("O")
$ws.Range("B8").Value = "O"
$ws.Range("B10").Value = "O"
$ws.Range("B11").Value = "O"

# Update the active selection to C8
$ws.Range("C8").Select()
